# fix bugs after changing column order
# Column E (the "list-type" marker, previously misaligned one column to the
# left of where it belongs) is swapped with column F, and a brand-new column
# H is populated with a copy of that same marker value, for every detail row
# (rows 7-40). Rows 2-6 only have a value in E (no F yet), so it is simply
# moved over into F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# ---------------------------------------------------------------------
# Step 1: populate new column H (rows 6-40) with the "목록/단수/OX" marker,
# using the same style (s=24) that column F already carries on these rows.
# Do this before the E/F swap below, while F still holds that style.
# ---------------------------------------------------------------------
$hValues = [ordered]@{
    6  = "목록"
    7  = "목록"
    8  = "단수"
    9  = "목록"
    10 = "목록"
    11 = "단수"
    12 = "목록"
    13 = "목록"
    14 = "단수"
    15 = "단수"
    16 = "목록"
    17 = "단수"
    18 = "단수"
    19 = "목록"
    20 = "목록"
    21 = "목록"
    22 = "단수"
    23 = "단수"
    24 = "OX"
    25 = "OX"
    26 = "단수"
    27 = "단수"
    28 = "OX"
    29 = "OX"
    30 = "단수"
    31 = "단수"
    32 = "목록"
    33 = "목록"
    34 = "OX"
    35 = "단수"
    36 = "OX"
    37 = "단수"
    38 = "OX"
    39 = "단수"
    # Row 40 intentionally has no H cell - it stays at E/F only.
}

$styleSource = $ws.Range("F7")
foreach ($r in $hValues.Keys) {
    $styleSource.Copy()
    $dst = $ws.Range("H$r")
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Value = $hValues[$r]
}

# ---------------------------------------------------------------------
# Step 2: rows 2-6 only use column E; move that value over to column F.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 6; $r++) {
    $e = $ws.Range("E$r")
    $f = $ws.Range("F$r")

    $e.Copy()
    $f.PasteSpecial($xlPasteFormats)
    $f.PasteSpecial($xlPasteValues)

    $e.Clear()
}

# ---------------------------------------------------------------------
# Step 3: rows 7-40 - swap columns E and F (value + format) via a
# scratch cell well outside the used range.
# ---------------------------------------------------------------------
$tmp = $ws.Range("Z1")
for ($r = 7; $r -le 40; $r++) {
    $e = $ws.Range("E$r")
    $f = $ws.Range("F$r")

    $e.Copy()
    $tmp.PasteSpecial($xlPasteFormats)
    $tmp.PasteSpecial($xlPasteValues)

    $f.Copy()
    $e.PasteSpecial($xlPasteFormats)
    $e.PasteSpecial($xlPasteValues)

    $tmp.Copy()
    $f.PasteSpecial($xlPasteFormats)
    $f.PasteSpecial($xlPasteValues)

    $tmp.Clear()
}

# ---------------------------------------------------------------------
# Step 4: move the active selection to F1, as in the saved workbook.
# ---------------------------------------------------------------------
$ws.Range("F1").Select()
